$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 11 new rows above the current row 2 (shifts old rows 2-21 down to 13-32)
$ws.Range("A2:A12").EntireRow.Insert()

# Newly inserted rows inherit the bold/centered header formatting - clear it so
# the new data cells have the same (default) style as the rest of the data rows.
$ws.Range("A2:C12").ClearFormats()

# The shift pushes the old last row (originally row 21) to row 32, beyond the
# new data range (A1:C31) - remove it so the sheet ends at row 31.
$ws.Rows.Item(32).Delete()

# Populate the newly inserted rows (2-12) with the new data values.
$newData = @(
  @(0.02334324724790524, 0.01564411001234516, 0.0184190768475939),
  @(-0.001020592069480415, 0.02263181184123195, -0.0006294894690920053),
  @(-0.02465064778197097, 0.04819875901065215, -0.01549884358920696),
  @(0.01120043709510704, 0.02745168796944902, 0.004529342418763651),
  @(-0.007870477419800848, 0.04201560953586564, 0.04054804218978412),
  @(-0.03585853518509281, 0.05588672146564573, 0.04165430539628354),
  @(0.01350235557410764, -0.04772198527324475, 0.04825462864303003),
  @(0.06278502832098701, -0.1541466276820112, 0.09075818756004657),
  @(0.02069492338270688, -0.3678749610738054, 0.1358281277665277),
  @(-0.1521724931350568, -0.4923760018697599, 0.1619388900878952),
  @(-0.2075079472326651, -0.5799234204175996, 0.3153703608890859)
)

$startRow = 2
for ($i = 0; $i -lt $newData.Count; $i++) {
    $rowIndex = $startRow + $i
    $values = $newData[$i]
    $ws.Cells.Item($rowIndex, 1).Value = $values[0]
    $ws.Cells.Item($rowIndex, 2).Value = $values[1]
    $ws.Cells.Item($rowIndex, 3).Value = $values[2]
}
